$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.147.42"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "'2.421.36"
$ws.Range("E3").Value = "  +0.09%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'554.15"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'137.43"
$ws.Range("E6").Value = "  -0.51%  "

$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("D8").Value = "'0.589"
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").Value = "'0.106"
$ws.Range("E9").Value = "  -1.21%  "

$ws.Range("D10").Value = "'5.71"
$ws.Range("E10").Value = "  -0.83%  "

$ws.Range("E11").Value = "  -0.24%  "

$ws.Range("D12").Value = "'0.353"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("D13").Value = "'24.92"
$ws.Range("E13").Value = "  +0.03%  "

$ws.Range("D14").Value = "'2.853.76"
$ws.Range("E14").Value = "  +0.02%  "

$ws.Range("D15").Value = "'60.073.56"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "'2.424.20"
$ws.Range("E17").Value = "  +0.83%  "

$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "'4.50"
$ws.Range("E19").Value = "  +2.44%  "

$ws.Range("D20").Value = "'327.34"
$ws.Range("E20").Value = "  -1.43%  "

$ws.Range("D21").Value = "'6.76"
$ws.Range("E21").Value = "  -0.25%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'65.14"
$ws.Range("E23").Value = "  +0.13%  "

$ws.Range("E24").Value = "  +4.12%  "

$ws.Range("D25").Value = "'8.75"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("E26").Value = "  -0.03%  "

$ws.Range("E27").Value = "  +3.35%  "

$ws.Range("D28").Value = "'0.0₃0774"
$ws.Range("E28").Value = "  -1.32%  "

$ws.Range("E29").Value = "  -0.43%  "

$ws.Range("D30").Value = "'170.03"
$ws.Range("E30").Value = "  +0.63%  "

$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("B32").Value = "PolygonEcosystemToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D32").Value = "'0.405"
$ws.Range("E32").Value = "  -2.65%  "

$ws.Range("B33").Value = "SuiNetwork"
$ws.Range("C33").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D33").Value = "'1.07"
$ws.Range("E33").Value = "  +2.56%  "

$ws.Range("D34").Value = "'18.55"
$ws.Range("E34").Value = "  -0.82%  "

$ws.Range("D36").Value = "'1.33"
$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("D38").Value = "'4.22"
$ws.Range("E38").Value = "  -0.17%  "

$ws.Range("D39").Value = "'328.43"
$ws.Range("E39").Value = "  +1.65%  "

$ws.Range("E40").Value = "  -0.53%  "

$ws.Range("D41").Value = "'144.97"
$ws.Range("E41").Value = "  +3.46%  "

$ws.Range("D42").Value = "'3.67"
$ws.Range("E42").Value = "  -0.61%  "

$ws.Range("D43").Value = "'20.04"
$ws.Range("E43").Value = "  +2.62%  "

$ws.Range("E44").Value = "  +0.45%  "

$ws.Range("D45").Value = "'0.0517"
$ws.Range("E45").Value = "  -0.67%  "

$ws.Range("D46").Value = "'0.576"
$ws.Range("E46").Value = "  +0.48%  "

$ws.Range("E47").Value = "  -1.35%  "

$ws.Range("D48").Value = "'11.04"
$ws.Range("E48").Value = "  -0.11%  "

$ws.Range("E49").Value = "  -2.16%  "

$ws.Range("E50").Value = "  -0.61%  "

$ws.Range("D51").Value = "'0.945"
$ws.Range("E51").Value = "  -0.61%  "

